$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "54.335.86"
$ws.Range("E2").Value = "  -7.70%  "
$ws.Range("D3").Value = "2.867.28"
$ws.Range("E3").Value = "  -10.61%  "
$ws.Range("D4").Value = "'1.00"
$ws.Range("E4").Value = "  +0.14%  "
$ws.Range("D5").Value = "'472.87"
$ws.Range("E5").Value = "  -11.51%  "
$ws.Range("D6").Value = "'126.09"
$ws.Range("E6").Value = "  -6.59%  "
$ws.Range("D7").Value = "'1.00"
$ws.Range("E7").Value = "  -0.04%  "
$ws.Range("D8").Value = "2.867.26"
$ws.Range("E8").Value = "  -10.53%  "
$ws.Range("D9").Value = "'0.402"
$ws.Range("E9").Value = "  -12.09%  "
$ws.Range("D10").Value = "'6.63"
$ws.Range("E10").Value = "  -11.75%  "
$ws.Range("D11").Value = "'0.0959"
$ws.Range("E11").Value = "  -15.78%  "
$ws.Range("D12").Value = "'0.331"
$ws.Range("E12").Value = "  -15.69%  "
$ws.Range("E13").Value = "  -4.68%  "
$ws.Range("D14").Value = "3.362.22"
$ws.Range("E14").Value = "  -10.58%  "
$ws.Range("D15").Value = "'23.05"
$ws.Range("E15").Value = "  -10.74%  "
$ws.Range("D16").Value = "54.379.17"
$ws.Range("E16").Value = "  -7.62%  "
$ws.Range("D17").Value = "2.868.41"
$ws.Range("E17").Value = "  -10.65%  "
$ws.Range("E18").Value = "  -14.67%  "
$ws.Range("D19").Value = "'5.32"
$ws.Range("E19").Value = "  -10.04%  "
$ws.Range("D20").Value = "'11.47"
$ws.Range("E20").Value = "  -13.51%  "
$ws.Range("D21").Value = "'7.07"
$ws.Range("E21").Value = "  -13.75%  "
$ws.Range("D22").Value = "'295.13"
$ws.Range("E22").Value = "  -18.18%  "
$ws.Range("E23").Value = "  +0.30%  "
$ws.Range("D24").Value = "'0.443"
$ws.Range("E24").Value = "  -14.49%  "
$ws.Range("D25").Value = "'58.68"
$ws.Range("E25").Value = "  -16.22%  "
$ws.Range("D26").Value = "'1.00"
$ws.Range("E26").Value = "  +0.00%  "
$ws.Range("E27").Value = "  -9.96%  "
$ws.Range("E28").Value = "  -0.08%  "
$ws.Range("D29").Value = "0.0₃0805"
$ws.Range("E29").Value = "  -16.47%  "
$ws.Range("D30").Value = "'6.24"
$ws.Range("E30").Value = "  -12.32%  "
$ws.Range("D31").Value = "'1.12"
$ws.Range("E31").Value = "  -7.07%  "
$ws.Range("D32").Value = "'6.19"
$ws.Range("E32").Value = "  -12.11%  "
$ws.Range("B33").Value = "PancakeSwap"
$ws.Range("C33").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D33").Value = "'1.61"
$ws.Range("E33").Value = "  -16.06%  "
$ws.Range("B34").Value = "EthereumClassic"
$ws.Range("C34").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D34").Value = "'18.90"
$ws.Range("E34").Value = "  -12.98%  "
$ws.Range("D35").Value = "'4.23"
$ws.Range("E35").Value = "  -13.67%  "
$ws.Range("D36").Value = "'135.00"
$ws.Range("E36").Value = "  -16.56%  "
$ws.Range("D37").Value = "'5.40"
$ws.Range("E37").Value = "  -15.11%  "
$ws.Range("E38").Value = "  -15.06%  "
$ws.Range("D39").Value = "'22.98"
$ws.Range("E39").Value = "  -11.85%  "
$ws.Range("D40").Value = "'0.0614"
$ws.Range("E40").Value = "  -13.00%  "
$ws.Range("D41").Value = "2.895.11"
$ws.Range("E41").Value = "  -10.54%  "
$ws.Range("D42").Value = "'0.999"
$ws.Range("E42").Value = "  -0.05%  "
$ws.Range("D43").Value = "'35.22"
$ws.Range("E43").Value = "  -13.79%  "
$ws.Range("D44").Value = "'0.959"
$ws.Range("E44").Value = "  -12.67%  "
$ws.Range("E45").Value = "  -15.84%  "
$ws.Range("E46").Value = "  -12.30%  "
$ws.Range("D47").Value = "'3.37"
$ws.Range("E47").Value = "  -16.01%  "
$ws.Range("D48").Value = "2.044.13"
$ws.Range("E48").Value = "  -11.25%  "
$ws.Range("D49").Value = "'5.30"
$ws.Range("E49").Value = "  -15.35%  "
$ws.Range("D50").Value = "'17.88"
$ws.Range("E50").Value = "  -13.80%  "
$ws.Range("D51").Value = "'0.0211"
$ws.Range("E51").Value = "  -11.78%  "
